$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 5 (duplicate fastq/rna sample entry), shifting all
# subsequent rows up by one.
$ws.Rows("5:5").Delete()

# Select the new row 5 (entire row), matching Excel's post-delete selection.
$ws.Rows("5:5").Select()
